$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1ecf9b45-3cb4-44f7-8406-405330071b14", "Check if the data is correct", "FAILED", "2023-10-18T14:51:44.079540800", "2023-10-18T14:51:49.382178900", "PT5.3026381S"),
    @("340a6953-2cd7-4db1-bfd6-65c534d290a4", "Check if the data is correct", "FAILED", "2023-10-18T14:52:45.036011200", "2023-10-18T14:52:49.303162200", "PT4.267151S"),
    @("adf5a47e-f644-4519-8965-263608b3af15", "Check if the data is correct", "FAILED", "2023-10-18T14:53:42.420576", "2023-10-18T14:53:47.147197900", "PT4.7266219S"),
    @("e68beea1-1508-45de-a585-6239a4beb37b", "Check if the data is correct", "PASSED", "2023-10-18T14:54:40.656899600", "2023-10-18T14:54:48.843654200", "PT8.1867546S")
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
